$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Remove now-unused formatted cells (outside the new C6:I19 extent) ---
$ws.Range("J6:P6").Clear()
$ws.Range("D7:H8").Clear()
$ws.Range("D12:H13").Clear()

# --- Widen column I and retire the old "bestFit" width ---
$ws.Columns("I").ColumnWidth = 24.7265625

# --- New test-case rows: With / Without spare in the 10th frame ---
$ws.Range("C17:C18").Merge()
$ws.Range("D17:D18").Merge()
$ws.Range("C17:D18").HorizontalAlignment = -4131
$ws.Range("C17:D18").VerticalAlignment = -4108

$ws.Range("I17").Value = "With Spare in 10th frame"
$ws.Range("E18").Value = "rolls = {2,4,3,6,4,5,7,2,3,5,10,0,1,7,3,4,5,9,0}"
$ws.Range("F18").Value = 85
$ws.Range("I18").Value = "Without Spare in 10th frame"

# --- Selection / view state ---
$ws.Range("D21").Select()
